$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume snapshot values (and two ranking swaps: SuiNetwork<->Dai, dogwifhat<->OKB).
$updates = @{
    'D2' = '61.190.91'
    'E2' = '  -0.06%  '
    'D3' = '2.404.69'
    'E4' = '  +0.50%  '
    'D5' = '568.00'
    'E5' = '  -0.39%  '
    'D6' = '142.05'
    'E6' = '  +1.45%  '
    'D7' = '0.999'
    'E8' = '  +1.49%  '
    'D9' = '2.413.96'
    'E9' = '  +0.05%  '
    'E10' = '  +1.45%  '
    'E11' = '  -0.15%  '
    'E12' = '  +2.59%  '
    'E13' = '  +2.14%  '
    'D14' = '26.45'
    'E14' = '  +1.07%  '
    'E15' = '  -0.27%  '
    'D16' = '2.814.50'
    'E16' = '  -0.95%  '
    'D17' = '60.833.08'
    'E17' = '  -0.43%  '
    'D18' = '2.471.51'
    'E18' = '  +2.62%  '
    'D19' = '8.15'
    'E19' = '  +4.24%  '
    'D20' = '10.69'
    'D21' = '324.00'
    'E21' = '  +0.11%  '
    'E22' = '  +0.68%  '
    'E23' = '  -0.40%  '
    'B24' = 'SuiNetwork'
    'C24' = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
    'D24' = '1.94'
    'E24' = '  +5.37%  '
    'B25' = 'Dai'
    'C25' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D25' = '1.00'
    'E25' = '  -0.18%  '
    'D26' = '65.19'
    'E26' = '  +0.86%  '
    'D27' = '593.76'
    'E27' = '  +0.44%  '
    'D28' = '8.27'
    'E28' = '  -0.06%  '
    'D29' = '0.0₃0948'
    'E29' = '  +1.22%  '
    'D30' = '2.522.27'
    'E30' = '  -0.72%  '
    'D31' = '8.02'
    'E31' = '  +1.34%  '
    'E32' = '  +1.58%  '
    'D33' = '1.81'
    'E33' = '  -1.08%  '
    'E34' = '  -0.74%  '
    'D35' = '0.997'
    'E35' = '  -0.49%  '
    'E36' = '  +2.56%  '
    'D37' = '0.373'
    'E37' = '  +1.13%  '
    'E38' = '  +0.00%  '
    'D39' = '152.16'
    'E39' = '  +0.19%  '
    'D40' = '18.36'
    'E40' = '  +0.65%  '
    'D41' = '5.27'
    'E41' = '  +1.68%  '
    'E42' = '  -0.01%  '
    'E43' = '  +1.11%  '
    'B44' = 'dogwifhat'
    'C44' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D44' = '2.52'
    'E44' = '  +6.41%  '
    'B45' = 'OKB'
    'C45' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D45' = '41.68'
    'E45' = '  +1.10%  '
    'E46' = '  -3.03%  '
    'D47' = '141.54'
    'E47' = '  -1.38%  '
    'D48' = '3.53'
    'E48' = '  -0.06%  '
    'D49' = '19.86'
    'E49' = '  +1.54%  '
    'E50' = '  +0.56%  '
    'E51' = '  +1.44%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "568.00") are not
    # auto-converted to real numbers, then restore the default "Normal" style
    # so no stray number-format style is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}

Write-Host "Applied cryptos update ($($updates.Count) cells)"
